$d = $word.ActiveDocument

# 1. Merge the two runs that were split around the old "_GoBack" bookmark
#    ("...best ha" + bookmark + "nd at the end...") back into a single run
#    with the full, uninterrupted text. This also removes the bookmark that
#    used to sit in the middle of "hand".
$d.Content.Find.Execute(
    "Players or AIs who has the best hand at the end of each round receives",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Players or AIs who has the best hand at the end of each round receives",
    2)

# 2. Re-create the "_GoBack" bookmark as a zero-length bookmark at the very
#    start of the document. A directly-collapsed Range at offset 0 can get
#    snapped to span the whole first paragraph, so instead insert a
#    temporary placeholder character, wrap the bookmark around it, then
#    delete the placeholder -- the bookmark collapses back down to a
#    zero-length bookmark exactly at the start, the way Word itself leaves
#    "_GoBack" after an edit.
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")

$bmRange = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range(0, 1).Delete()
